$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "TextBox 24" (id=25): featureTableJSON -> "JSON Header", reposition/resize ---
$shJson = $s.Shapes.Item("TextBox 24")
$shJson.Left = 51.7184
$shJson.Width = 81.1852
$shJson.TextFrame.TextRange.Text = "JSON Header"

# --- "TextBox 6" (id=7): drop the stray endParaRPr on the "POSITION:" paragraph ---
$shPos = $s.Shapes.Item("TextBox 6")
$trPos = $shPos.TextFrame.TextRange
$firstPara = $trPos.Paragraphs(1, 1)
[void]$firstPara.Delete()
$trPosAfter = $shPos.TextFrame.TextRange
[void]$trPosAfter.InsertBefore("POSITION:" + [char]13)

# --- "TextBox 14" (id=15): reposition only (text/formatting unchanged) ---
$shLen = $s.Shapes.Item("TextBox 14")
$shLen.Left = 332.899
$shLen.Top = 83.3307

# --- "TextBox 19" (id=20): featureTableBinary -> "Binary Body", reposition/resize ---
$shBin = $s.Shapes.Item("TextBox 19")
$shBin.Left = 257.4936
$shBin.Width = 81.1852
$shBin.TextFrame.TextRange.Text = "Binary Body"
